# Updates the inventory report: refresh the timestamp, rename the "ID"
# header to "N°.", restate rows 3-4 with new sample data, and remove the
# now-unused rows 5-13 (shrinking the sheet from A1:F13 down to A1:F4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Refresh the generated-at timestamp in the title row.
$ws.Range("A1").Value = "Reporte generado el 19/05/2025 a las 11:35"

# 2. Rename the "ID" column header to "N°."
$ws.Range("A2").Value = "N°."

# 3. Row 3: new sample item.
$ws.Range("B3").Value = "Mueble 123"
$ws.Range("C3").Value = 123
$ws.Range("D3").Value = 321
$ws.Range("E3").Value = "juego"
$ws.Range("F3").Value = "Sala"

# 4. Row 4: new sample item.
$ws.Range("B4").Value = "Boe 123"
$ws.Range("C4").Value = 43
$ws.Range("D4").Value = 35
$ws.Range("F4").Value = "Dormitorio"

# 5. Remove rows 5 through 13 entirely - only 4 data rows remain now.
$ws.Range("A5:F13").EntireRow.Delete()
